# Apply updated 'CPU time taken for sorting' values (column C) to rows 2-101
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 131717.0
    3 = 297218.0
    4 = 354096.0
    5 = 477261.0
    6 = 563646.0
    7 = 1061434.0
    8 = 793724.0
    9 = 929717.0
    10 = 1028078.0
    11 = 1162788.0
    12 = 1286380.0
    13 = 1405695.0
    14 = 1525438.0
    15 = 1643898.0
    16 = 3046599.0
    17 = 1895358.0
    18 = 2067701.0
    19 = 2174187.0
    20 = 2244750.0
    21 = 2406831.0
    22 = 3386156.0
    23 = 2600557.0
    24 = 3851869.0
    25 = 2995281.0
    26 = 3026927.0
    27 = 3156934.0
    28 = 3233911.0
    29 = 3528992.0
    30 = 3516162.0
    31 = 3559355.0
    32 = 5101044.0
    33 = 3879667.0
    34 = 4035333.0
    35 = 4120863.0
    36 = 4286792.0
    37 = 4467689.0
    38 = 4488645.0
    39 = 4738394.0
    40 = 4697767.0
    41 = 4950082.0
    42 = 4873960.0
    43 = 5133546.0
    44 = 5439745.0
    45 = 5252005.0
    46 = 5426488.0
    47 = 5548369.0
    48 = 5810092.0
    49 = 6833893.0
    50 = 5861410.0
    51 = 6082934.0
    52 = 6451572.0
    53 = 7542514.0
    54 = 7729399.0
    55 = 6500752.0
    56 = 8210936.0
    57 = 6769317.0
    58 = 8892187.0
    59 = 6990414.0
    60 = 7434318.0
    61 = 7254704.0
    62 = 7640019.0
    63 = 9035451.0
    64 = 7896611.0
    65 = 9014495.0
    66 = 8015498.0
    67 = 8024906.0
    68 = 8319559.0
    69 = 8816920.0
    70 = 8836592.0
    71 = 9523830.0
    72 = 8676222.0
    73 = 10694315.0
    74 = 8924688.0
    75 = 9055550.0
    76 = 9564030.0
    77 = 9318556.0
    78 = 9415634.0
    79 = 10175146.0
    80 = 11221612.0
    81 = 9788975.0
    82 = 9908718.0
    83 = 10062673.0
    84 = 12281764.0
    85 = 12269789.0
    86 = 12363017.0
    87 = 10647274.0
    88 = 12673066.0
    89 = 10829882.0
    90 = 11263094.0
    91 = 11307998.0
    92 = 11328953.0
    93 = 11468796.0
    94 = 11488040.0
    95 = 11592387.0
    96 = 12176133.0
    97 = 11826740.0
    98 = 11952043.0
    99 = 12718825.0
    100 = 14414896.0
    101 = 12368578.0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
